$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 537
$ws.Range("I15").Value = 537
$ws.Range("K15").Value = 1611
$ws.Range("M15").Value = -1442
$ws.Range("H28").Value = 529.7273
$ws.Range("I28").Value = 382.7
$ws.Range("K28").Value = 382.7
$ws.Range("M28").Value = 102.3
$ws.Range("H39").Value = 88.72727
$ws.Range("I39").Value = 17.6
$ws.Range("J39").Value = 800
$ws.Range("K39").Value = 52.8
$ws.Range("L39").Value = 2400
$ws.Range("M39").Value = 243.2
$ws.Range("N39").Value = -2992
$ws.Range("H41").Value = 578.625
$ws.Range("J41").Value = 656.8570999999999
$ws.Range("L41").Value = 656.8570999999999
$ws.Range("N41").Value = -1536.8571
$ws.Range("H53").Value = 866.25
$ws.Range("J53").Value = 591.8
$ws.Range("L53").Value = 591.8
$ws.Range("N53").Value = -1865.8
$ws.Range("H62").Value = 4636
$ws.Range("I62").Value = 3654.5
$ws.Range("J62").Value = 5028.6
$ws.Range("K62").Value = 3654.5
$ws.Range("L62").Value = 5028.6
$ws.Range("M62").Value = -3030.5
$ws.Range("N62").Value = -6276.6
$ws.Range("H65").Value = 4636
$ws.Range("I65").Value = 3654.5
$ws.Range("J65").Value = 5028.6
$ws.Range("K65").Value = 18272.5
$ws.Range("L65").Value = 25143
$ws.Range("M65").Value = -15152.5
$ws.Range("N65").Value = -31383
$ws.Range("H86").Value = 7662.75
$ws.Range("I86").Value = 6799.6665
$ws.Range("K86").Value = 6799.6665
$ws.Range("M86").Value = -5676.6665
$ws.Range("H89").Value = 7662.75
$ws.Range("I89").Value = 6799.6665
$ws.Range("K89").Value = 33998.3325
$ws.Range("M89").Value = -28382.3325
$ws.Range("H92").Value = 624.3333
$ws.Range("I92").Value = 588.2
$ws.Range("J92").Value = 805
$ws.Range("K92").Value = 588.2
$ws.Range("L92").Value = 805
$ws.Range("M92").Value = 659.8
$ws.Range("N92").Value = -3301
$ws.Range("H98").Value = 1999
$ws.Range("I98").Value = 1999
$ws.Range("K98").Value = 1999
$ws.Range("M98").Value = -501
$ws.Range("H106").Value = 1999.5
$ws.Range("I106").Value = 1999.5
$ws.Range("K106").Value = 1999.5
$ws.Range("M106").Value = -1368.5
$ws.Range("H107").Value = 336.7143
$ws.Range("I107").Value = 289.92307
$ws.Range("K107").Value = 289.92307
$ws.Range("M107").Value = 1630.07693
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1999
$ws.Range("K122").Value = 5997
$ws.Range("M122").Value = -3547
$ws.Range("H125").Value = 83337280
$ws.Range("I125").Value = 250000000
$ws.Range("J125").Value = 5918
$ws.Range("K125").Value = 2250000000
$ws.Range("L125").Value = 53262
$ws.Range("M125").Value = -2249997540
$ws.Range("N125").Value = -58182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 11
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 11
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 339
$ws.Range("N15").ClearContents()
$ws.Range("H61").Value = 8255.5
$ws.Range("I61").Value = 19999
$ws.Range("K61").Value = 19999
$ws.Range("M61").Value = -19787
$ws.Range("H110").Value = 2548.5
$ws.Range("I110").Value = 3000
$ws.Range("J110").Value = 1194
$ws.Range("K110").Value = 3000
$ws.Range("L110").Value = 1194
$ws.Range("M110").Value = -955
$ws.Range("N110").Value = -5284
$ws.Range("H136").Value = 8255.5
$ws.Range("I136").Value = 19999
$ws.Range("K136").Value = 59997
$ws.Range("M136").Value = -57447

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 3000
$ws.Range("I36").Value = 3000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2466

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H86").Value = 13944394
$ws.Range("J86").Value = 4500
$ws.Range("L86").Value = 4500
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 13944394
$ws.Range("J89").Value = 4500
$ws.Range("L89").Value = 22500
$ws.Range("N89").Value = -33732
$ws.Range("H98").Value = 70000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 70000
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 70000
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -74492
$ws.Range("H99").Value = 5945.5884
$ws.Range("I99").Value = 5312.769
$ws.Range("K99").Value = 5312.769
$ws.Range("M99").Value = -3814.769
$ws.Range("H106").Value = 17835.5
$ws.Range("J106").Value = 17835.5
$ws.Range("L106").Value = 17835.5
$ws.Range("N106").Value = -20359.5
$ws.Range("H107").Value = 500.9091
$ws.Range("I107").Value = 224.70589
$ws.Range("J107").Value = 1440
$ws.Range("K107").Value = 224.70589
$ws.Range("L107").Value = 1440
$ws.Range("M107").Value = 1695.29411
$ws.Range("N107").Value = -5280
$ws.Range("H122").Value = 980.625
$ws.Range("I122").Value = 795
$ws.Range("J122").Value = 1125
$ws.Range("K122").Value = 2385
$ws.Range("L122").Value = 3375
$ws.Range("M122").Value = 65
$ws.Range("N122").Value = -8275
$ws.Range("H126").Value = 5945.5884
$ws.Range("I126").Value = 5312.769
$ws.Range("K126").Value = 15938.307
$ws.Range("M126").Value = -13468.307
$ws.Range("H134").Value = 2421.9412
$ws.Range("I134").Value = 2012.3572
$ws.Range("J134").Value = 4333.3335
$ws.Range("K134").Value = 6037.071599999999
$ws.Range("L134").Value = 13000.0005
$ws.Range("M134").Value = -3502.071599999999
$ws.Range("N134").Value = -18070.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 898.3333
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 898.3333
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 2694.9999
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -3142.9999
$ws.Range("H12").Value = 1168.6666
$ws.Range("J12").Value = 1394
$ws.Range("L12").Value = 4182
$ws.Range("N12").Value = -4528
$ws.Range("H125").Value = 15000
$ws.Range("J125").Value = 15000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("H126").Value = 300
$ws.Range("I126").Value = 300
$ws.Range("K126").Value = 900
$ws.Range("M126").Value = 4040
$ws.Range("H131").Value = 1344.2
$ws.Range("I131").Value = 850
$ws.Range("J131").Value = 1556
$ws.Range("K131").Value = 2550
$ws.Range("L131").Value = 4668
$ws.Range("M131").Value = 2490
$ws.Range("N131").Value = -14748
$ws.Range("H132").Value = 2149.2856
$ws.Range("I132").Value = 1249.5
$ws.Range("J132").Value = 2509.2
$ws.Range("K132").Value = 11245.5
$ws.Range("L132").Value = 22582.8
$ws.Range("M132").Value = -8715.5
$ws.Range("N132").Value = -27642.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 116.791664
$ws.Range("I2").Value = 100.42105
$ws.Range("J2").Value = 179
$ws.Range("K2").Value = 100.42105
$ws.Range("L2").Value = 179
$ws.Range("M2").Value = 12.57895000000001
$ws.Range("N2").Value = -405
$ws.Range("H102").Value = 899
$ws.Range("I102").Value = 899
$ws.Range("K102").Value = 899
$ws.Range("M102").Value = 723
$ws.Range("H122").Value = 1882.3334
$ws.Range("I122").Value = 1882.3334
$ws.Range("K122").Value = 5647.0002
$ws.Range("M122").Value = -3197.0002
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 250005730
$ws.Range("I126").Value = 250005730
$ws.Range("K126").Value = 750017190
$ws.Range("M126").Value = -750014720
$ws.Range("H132").Value = 3946.3125
$ws.Range("I132").Value = 3946
$ws.Range("K132").Value = 11838
$ws.Range("M132").Value = -9308

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4991.6665
$ws.Range("I7").Value = 4991.6665
$ws.Range("K7").Value = 4991.6665
$ws.Range("M7").Value = -4879.6665
$ws.Range("H36").Value = 39999
$ws.Range("J36").Value = 39999
$ws.Range("L36").Value = 39999
$ws.Range("N36").Value = -41123
$ws.Range("H40").Value = 3323.1667
$ws.Range("I40").Value = 3323.1667
$ws.Range("K40").Value = 3323.1667
$ws.Range("M40").Value = -3187.1667
$ws.Range("H126").Value = 4991.6665
$ws.Range("I126").Value = 4991.6665
$ws.Range("K126").Value = 14974.9995
$ws.Range("M126").Value = -12504.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 24995
$ws.Range("I70").Value = 24990
$ws.Range("K70").Value = 24990
$ws.Range("M70").Value = -24675
$ws.Range("H73").Value = 24995
$ws.Range("I73").Value = 24990
$ws.Range("K73").Value = 24990
$ws.Range("M73").Value = -23898
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("H126").Value = 1936.8
$ws.Range("I126").Value = 1922.25
$ws.Range("K126").Value = 5766.75
$ws.Range("M126").Value = -3296.75
$ws.Range("H136").Value = 2385.8235
$ws.Range("I136").Value = 1859.2222
$ws.Range("J136").Value = 4417
$ws.Range("K136").Value = 5577.6666
$ws.Range("L136").Value = 13251
$ws.Range("M136").Value = -3027.6666
$ws.Range("N136").Value = -18351
